$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-01-29 Wednesday"; new="2025-01-30 Thursday"},
    @{old="435÷4="; new="832÷2="},
    @{old="384÷6="; new="890÷2="},
    @{old="764÷2="; new="845÷7="},
    @{old="282÷2="; new="870÷6="},
    @{old="318÷4="; new="807÷8="},
    @{old="658÷5="; new="113÷7="},
    @{old="491÷6="; new="366÷4="},
    @{old="237÷7="; new="556÷3="},
    @{old="589÷6="; new="708÷6="},
    @{old="130÷5="; new="319÷5="},
    @{old="811÷7="; new="164÷2="},
    @{old="249÷5="; new="504÷9="},
    @{old="406÷3="; new="942÷3="},
    @{old="530÷7="; new="349÷5="},
    @{old="732÷3="; new="131÷4="},
    @{old="819÷6="; new="611÷9="},
    @{old="574÷8="; new="437÷4="},
    @{old="316÷3="; new="982÷9="},
    @{old="968÷3="; new="320÷2="},
    @{old="498÷6="; new="400÷7="},
    @{old="462÷8="; new="779÷7="},
    @{old="933÷9="; new="396÷5="},
    @{old="654÷4="; new="991÷2="},
    @{old="349÷7="; new="938÷6="},
    @{old="604÷5="; new="705÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}

$d.Save()
